$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string (A1)
$ws.Range("A1").Value = "Datos actualizados a 7 de Octubre de 2020 a las 13:39"

# Swap country names that moved position in the source data feed
$ws.Range("A93").Value = "Senegal"
$ws.Range("A94").Value = "Zambia"
$ws.Range("A207").Value = "Nueva Caledonia"
$ws.Range("A208").Value = "Santa Lucia"

# Update numeric statistics for affected rows
# Row 16
$ws.Range("B16").Value = 483844
$ws.Range("C16").Value = 4019
$ws.Range("D16").Value = 397109
$ws.Range("E16").Value = 59077
$ws.Range("G16").Value = 239
$ws.Range("H16").Value = 27658

# Row 26
$ws.Range("B26").Value = 307456
$ws.Range("C26").Value = 337
$ws.Range("E26").Value = 32221

# Row 36
$ws.Range("B36").Value = 127181
$ws.Range("C36").Value = 238
$ws.Range("D36").Value = 124108
$ws.Range("E36").Value = 2855
$ws.Range("G36").Value = 2
$ws.Range("H36").Value = 218

# Row 47
$ws.Range("B47").Value = 94253
$ws.Range("C47").Value = 3439
$ws.Range("D47").Value = 68668
$ws.Range("E47").Value = 25007
$ws.Range("G47").Value = 15
$ws.Range("H47").Value = 578

# Row 62
$ws.Range("B62").Value = 57709
$ws.Range("C62").Value = 1077
$ws.Range("E62").Value = 8327
$ws.Range("G62").Value = 1
$ws.Range("H62").Value = 2082

# Row 92
$ws.Range("B92").Value = 16633
$ws.Range("C92").Value = 33
$ws.Range("D92").Value = 15808
$ws.Range("E92").Value = 590
$ws.Range("G92").Value = 1
$ws.Range("H92").Value = 235

# Row 93
$ws.Range("B93").Value = 15174
$ws.Range("C93").Value = 33
$ws.Range("D93").Value = 12998
$ws.Range("E93").Value = 1863
$ws.Range("G93").Value = 1
$ws.Range("H93").Value = 313

# Row 94
$ws.Range("B94").Value = 15170
$ws.Range("D94").Value = 14313
$ws.Range("E94").Value = 522
$ws.Range("H94").Value = 335

# Row 108
$ws.Range("E108").Value = 4455
$ws.Range("G108").Value = 5
$ws.Range("H108").Value = 63

# Row 144
$ws.Range("B144").Value = 3442
$ws.Range("C144").Value = 68
$ws.Range("D144").Value = 2865
$ws.Range("E144").Value = 536
$ws.Range("G144").Value = 1
$ws.Range("H144").Value = 41

# Row 168
$ws.Range("B168").Value = 1099
$ws.Range("C168").Value = 1
$ws.Range("E168").Value = 41

# Row 175
$ws.Range("B175").Value = 523
$ws.Range("C175").Value = 2
$ws.Range("E175").Value = 31

